$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before what is currently row 1250, pushing the
# existing data (old rows 1250-1314) down to rows 1252-1316.
$ws.Rows.Item(1250).Resize(2).Insert()

# New row 1250: Acelga, Primera, week of date 45147
$ws.Cells.Item(1250,1).Value = 6
$ws.Cells.Item(1250,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1250,3).Value = "Metropolitana"
$ws.Cells.Item(1250,4).Value = 45147
$ws.Cells.Item(1250,5).Value = 13
$ws.Cells.Item(1250,6).Value = 100112009
$ws.Cells.Item(1250,7).Value = "Acelga"
$ws.Cells.Item(1250,8).Value = "Sin especificar"
$ws.Cells.Item(1250,9).Value = "Primera"
$ws.Cells.Item(1250,10).Value = 250
$ws.Cells.Item(1250,11).Value = 10000
$ws.Cells.Item(1250,12).Value = 10000
$ws.Cells.Item(1250,13).Value = 10000
$ws.Cells.Item(1250,14).Value = "`$/docena de atados"
$ws.Cells.Item(1250,15).Value = "Región Metropolitana"
$ws.Cells.Item(1250,16).Value = 3333
$ws.Cells.Item(1250,17).Value = 3
$ws.Cells.Item(1250,18).Value = "Hortaliza"

# New row 1251: Acelga, Segunda, week of date 45147
$ws.Cells.Item(1251,1).Value = 6
$ws.Cells.Item(1251,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1251,3).Value = "Metropolitana"
$ws.Cells.Item(1251,4).Value = 45147
$ws.Cells.Item(1251,5).Value = 13
$ws.Cells.Item(1251,6).Value = 100112009
$ws.Cells.Item(1251,7).Value = "Acelga"
$ws.Cells.Item(1251,8).Value = "Sin especificar"
$ws.Cells.Item(1251,9).Value = "Segunda"
$ws.Cells.Item(1251,10).Value = 170
$ws.Cells.Item(1251,11).Value = 8000
$ws.Cells.Item(1251,12).Value = 8000
$ws.Cells.Item(1251,13).Value = 8000
$ws.Cells.Item(1251,14).Value = "`$/docena de atados"
$ws.Cells.Item(1251,15).Value = "Región Metropolitana"
$ws.Cells.Item(1251,16).Value = 2667
$ws.Cells.Item(1251,17).Value = 3
$ws.Cells.Item(1251,18).Value = "Hortaliza"
